$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 5; $r -le 16; $r++) {
    $ws.Cells.Item($r, 1).ClearContents()
    $ws.Cells.Item($r, 2).ClearContents()
}

$ws.Range("H14").Select()
